$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Startup")

# Add the two new startups discovered - Xanadu.ai first (location filled before CEO name,
# which was tracked down and added afterward), then ANELLO Photonics.
$ws.Range("B34").Value = "Xanadu.ai"
$ws.Range("C34").Value = 190
$ws.Range("D34").Value = "Quantum Computing / Optical Comptuing, caling to one million qubits through optical networking"
$ws.Range("F34").Value = "Toronto, Canada"

$ws.Range("B35").Value = "ANELLO Photonics"
$ws.Range("C35").Value = 30
$ws.Range("D35").Value = "Silicon Photonics for Fiber Optical Gyroscope tech (FOGs)"
$ws.Range("E35").Value = "CEO: Mario Paniccia and CTO Mike Horton"
$ws.Range("F35").Value = "Santa Clara, CA"

$ws.Range("E34").Value = "Christian Weedbrook"

# Bring the Startup sheet back into focus (it was left scrolled to row 17 on the
# Defense Contractors tab before); select G6 there.
$ws.Activate()
[void]$ws.Range("G6").Select()
